$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "Hyperlink" named style (currently applied to D2, the product-URL
# cell) before we touch anything, so we can re-apply the exact same style object
# to the cell the hyperlink ends up in after the column shift.
$hyperlinkStyle = $ws.Range("D2").Style

# Delete column A (the "Id" column) entirely; B:E shift left to become A:D.
$ws.Columns("A").Delete()

# The hyperlink that used to anchor on D2 now sits on the shifted-left cell C2,
# but this runtime does not auto-re-anchor Hyperlink objects when a column is
# deleted, so recreate it explicitly pointing at the new location.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "http://www.prod1.com/")

# Re-apply the original Hyperlink cell style (Add() re-derives its own), so C2
# ends up referencing the same style as before rather than a new duplicate one.
$ws.Range("C2").Style = $hyperlinkStyle

# Update the active selection to match the edited workbook's saved state.
$ws.Range("C11").Select()
